$wb = $excel.ActiveWorkbook

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4642.6
$ws.Range("I100").Value = 3572.875
$ws.Range("K100").Value = 3572.875
$ws.Range("M100").Value = -3031.875

# ALC row 105
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 54720.832
$ws.Range("J105").Value = 54720.832
$ws.Range("L105").Value = 54720.832
$ws.Range("N105").Value = -61708.832

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 8005.8887
$ws.Range("I131").Value = 7144.6665
$ws.Range("K131").Value = 21433.9995
$ws.Range("M131").Value = -16393.9995

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 36799332
$ws.Range("I5").Value = 36799332
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 36799332
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -36799220

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1417.5
$ws.Range("I74").Value = 1533
$ws.Range("K74").Value = 1533
$ws.Range("M74").Value = -659

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1417.5
$ws.Range("I77").Value = 1533
$ws.Range("K77").Value = 7665
$ws.Range("M77").Value = -3297

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 142861070
$ws.Range("I132").Value = 250004180
$ws.Range("K132").Value = 750012540
$ws.Range("M132").Value = -750010010

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 36799332
$ws.Range("I4").Value = 36799332
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 36799332
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -36799217

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 537.5
$ws.Range("I22").Value = 416.66666
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 416.66666
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -243.66666
$ws.Range("N22").Value = -1246

# BSM row 27
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 32500
$ws.Range("J27").Value = 32500
$ws.Range("L27").Value = 32500
$ws.Range("N27").Value = -32884

# BSM row 59
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 97500
$ws.Range("J59").Value = 97500
$ws.Range("L59").Value = 97500
$ws.Range("N59").Value = -99194

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 518.75
$ws.Range("I80").Value = 417
$ws.Range("J80").Value = 552.6667
$ws.Range("K80").Value = 417
$ws.Range("L80").Value = 552.6667
$ws.Range("M80").Value = 581
$ws.Range("N80").Value = -2548.6667

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 518.75
$ws.Range("I83").Value = 417
$ws.Range("J83").Value = 552.6667
$ws.Range("K83").Value = 2085
$ws.Range("L83").Value = 2763.3335
$ws.Range("M83").Value = 2907
$ws.Range("N83").Value = -12747.3335

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1350
$ws.Range("I134").Value = 1350
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4050
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -1515

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3591.5833
$ws.Range("I62").Value = 3771.4285
$ws.Range("J62").Value = 3339.8
$ws.Range("K62").Value = 3771.4285
$ws.Range("L62").Value = 3339.8
$ws.Range("M62").Value = -3147.4285
$ws.Range("N62").Value = -4587.8

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3591.5833
$ws.Range("I65").Value = 3771.4285
$ws.Range("J65").Value = 3339.8
$ws.Range("K65").Value = 18857.1425
$ws.Range("L65").Value = 16699
$ws.Range("M65").Value = -15737.1425
$ws.Range("N65").Value = -22939

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 35901440
$ws.Range("J99").Value = 112501250
$ws.Range("L99").Value = 112501250
$ws.Range("N99").Value = -112504246

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 31572418
$ws.Range("I122").Value = 50512260
$ws.Range("K122").Value = 151536780
$ws.Range("M122").Value = -151534330

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 35901440
$ws.Range("J126").Value = 112501250
$ws.Range("L126").Value = 337503750
$ws.Range("N126").Value = -337508690

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2035.8286
$ws.Range("I132").Value = 1613.8182
$ws.Range("J132").Value = 8999
$ws.Range("K132").Value = 4841.4546
$ws.Range("L132").Value = 26997
$ws.Range("M132").Value = -2311.4546
$ws.Range("N132").Value = -32057

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 816.3570999999999
$ws.Range("I12").Value = 728.3333
$ws.Range("J12").Value = 840.36365
$ws.Range("K12").Value = 2184.9999
$ws.Range("L12").Value = 2521.09095
$ws.Range("M12").Value = -2011.9999
$ws.Range("N12").Value = -2867.09095

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 1985
$ws.Range("I87").Value = 1985
$ws.Range("K87").Value = 5955
$ws.Range("M87").Value = -4707

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 1985
$ws.Range("I90").Value = 1985
$ws.Range("K90").Value = 17865
$ws.Range("M90").Value = -11625

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6030
$ws.Range("I122").Value = 7749.5
$ws.Range("K122").Value = 23248.5
$ws.Range("M122").Value = -20798.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8518.75
$ws.Range("I132").Value = 8017.35
$ws.Range("K132").Value = 24052.05
$ws.Range("M132").Value = -21522.05

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 38465120
$ws.Range("I122").Value = 58825356
$ws.Range("J122").Value = 6897.8887
$ws.Range("K122").Value = 176476068
$ws.Range("L122").Value = 20693.6661
$ws.Range("M122").Value = -176473618
$ws.Range("N122").Value = -25593.6661

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2974.182
$ws.Range("I132").Value = 2522.611
$ws.Range("K132").Value = 7567.833
$ws.Range("M132").Value = -5037.833

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31040

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2424.5715
$ws.Range("I122").Value = 2471.125
$ws.Range("J122").Value = 2362.5
$ws.Range("K122").Value = 7413.375
$ws.Range("L122").Value = 7087.5
$ws.Range("M122").Value = -4963.375
$ws.Range("N122").Value = -11987.5

Write-Output "Applied all changes"